# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The sheet holds a small metadata table (header + 3 metadata rows + an
# old "mapping file" row). The dimension/measure curation changed:
#  - "estado-civil" is now a measure (iaest-measure:...) instead of a
#    dimension (iaest-dimension:...)
#  - municipio-nombre / provincia-nombre / aragon / comarca-nombre are now
#    all curated as sdmx-dimension:refArea dimensions, each pointing to its
#    own URI mapping column (URI-Municipio / URI-Provincia / URI-Comunidad /
#    URI-comarca) instead of referencing external mapping-*.xlsx files.
#  - the old row 5 (mapping-estado-civil.xlsx / mapping-aragon.xlsx) is no
#    longer needed and is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-* dimension/measure identifiers
$ws.Range("B2").Value = "iaest-measure:estado-civil"
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "sdmx-dimension:refArea"

# Row 3: medida/dim role
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("D3").Value = "dim"
$ws.Range("E3").Value = "dim"
$ws.Range("F3").Value = "dim"
$ws.Range("H3").Value = "dim"

# Row 4: type / URI mapping columns
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("D4").Value = "URI-Municipio"
$ws.Range("E4").Value = "URI-Provincia"
$ws.Range("F4").Value = "URI-Comunidad"
$ws.Range("H4").Value = "URI-comarca"

# Row 5 (old external mapping-file references) is no longer used.
$ws.Rows.Item(5).Delete()
